# Apply edits to rules.xlsx / Sheet1:
#  - Repurpose existing column B ("GeneratePseudonym") into "SimplifiedPseudonym"
#    (values updated, column position/width untouched).
#  - Insert two brand-new columns C ("StandardPseudonym") and D ("RemoveField")
#    right after column B, pushing the former RequiredField/ValidateField/
#    ValidationRegExp columns from C,D,E to E,F,G.
#  - Add a new row 6 for field "Spouse" (RemoveField = Yes, everything else No).
#  - Resulting table occupies A1:G6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two blank columns before the old column C (RequiredField), i.e. at C:D.
# This shifts old C,D,E (RequiredField, ValidateField, ValidationRegExp) to E,F,G
# while leaving columns A and B completely untouched (no Insert/Delete applied
# to them), which preserves their original bestFit widths exactly.
$ws.Range("C:D").Insert() | Out-Null

# Give the two new columns (C,D) a width of 20, matching the non-bestFit
# "20" width used elsewhere in the sheet (ColumnWidth value chosen so the
# stored OOXML width comes out to exactly 20).
$ws.Range("C:D").ColumnWidth = 19.166666666666668

# --- Header row ---
# (StandardPseudonym is written before SimplifiedPseudonym so the shared-string
# table registers them in the same order as the canonical workbook.)
$ws.Range("A1").Value = "FieldName"
$ws.Range("C1").Value = "StandardPseudonym"
$ws.Range("B1").Value = "SimplifiedPseudonym"
$ws.Range("D1").Value = "RemoveField"
$ws.Range("E1").Value = "RequiredField"
$ws.Range("F1").Value = "ValidateField"
$ws.Range("G1").Value = "ValidationRegExp"

# --- Data rows ---
$ws.Range("A2").Value = "Name"
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "No"
$ws.Range("D2").Value = "No"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "No"

$ws.Range("A3").Value = "Age"
$ws.Range("B3").Value = "No"
$ws.Range("C3").Value = "No"
$ws.Range("D3").Value = "No"
$ws.Range("E3").Value = "Yes"
$ws.Range("F3").Value = "Yes"
$ws.Range("G3").Value = "[0-9]{1,3}"

$ws.Range("A4").Value = "Height-inches"
$ws.Range("B4").Value = "No"
$ws.Range("C4").Value = "No"
$ws.Range("D4").Value = "No"
$ws.Range("E4").Value = "Yes"
$ws.Range("F4").Value = "Yes"
$ws.Range("G4").Value = "[0-9]{1,3}"

$ws.Range("A5").Value = "Weight-pounds"
$ws.Range("B5").Value = "No"
$ws.Range("C5").Value = "No"
$ws.Range("D5").Value = "No"
$ws.Range("E5").Value = "Yes"
$ws.Range("F5").Value = "Yes"
$ws.Range("G5").Value = "[0-9]{1,3}"

$ws.Range("A6").Value = "Spouse"
$ws.Range("B6").Value = "No"
$ws.Range("C6").Value = "No"
$ws.Range("D6").Value = "Yes"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "No"

# Move the active selection to the next empty row below the table, matching
# the post-edit cursor position.
$ws.Range("A7").Select() | Out-Null
